$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-22 10:21:49"
$wsZhCn.Range("E5").Value = "2016-03-22 10:21:49"
$wsZhCn.Range("H3").Value = "2016-03-22 10:22:13"
$wsZhCn.Range("H5").Value = "2016-03-22 10:22:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-22 10:21:53"
$wsDeDe.Range("E5").Value = "2016-03-22 10:21:53"
$wsDeDe.Range("H3").Value = "2016-03-22 10:22:20"
$wsDeDe.Range("H5").Value = "2016-03-22 10:22:20"
